$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.281766
$ws.Range("H2").Value = 45.845298
$ws.Range("I2").Value = 0.1817381432449346
$ws.Range("J2").Value = 0.1817381432449346
$ws.Range("M2").Value = 1.910418
$ws.Range("N2").Value = 5.731254
$ws.Range("O2").Value = 0.01809124304049503
$ws.Range("P2").Value = 0.01809124304049503
$ws.Range("Q2").Value = 29.194560838188
$ws.Range("R2").Value = 262.751047543692
$ws.Range("S2").Value = 0.003287868919172412
$ws.Range("T2").Value = 0.003287868919172412
$ws.Range("G3").Value = 15.281766
$ws.Range("H3").Value = 45.845298
$ws.Range("I3").Value = 0.1817381432449346
$ws.Range("J3").Value = 0.1817381432449346
$ws.Range("O3").Value = 0.302988173785169
$ws.Range("P3").Value = 0.302988173785169
$ws.Range("Q3").Value = 488.944106992692
$ws.Range("R3").Value = 4400.496962934229
$ws.Range("S3").Value = 0.05506450812889019
$ws.Range("T3").Value = 0.05506450812889019
$ws.Range("G4").Value = 15.281766
$ws.Range("H4").Value = 45.845298
$ws.Range("I4").Value = 0.1817381432449346
$ws.Range("J4").Value = 0.1817381432449346
$ws.Range("M4").Value = 37.858701
$ws.Range("N4").Value = 113.576103
$ws.Range("O4").Value = 0.3585136661130873
$ws.Range("P4").Value = 0.3585136661130873
$ws.Range("Q4").Value = 578.5478097459659
$ws.Range("R4").Value = 5206.930287713693
$ws.Range("S4").Value = 0.06515560800732692
$ws.Range("T4").Value = 0.06515560800732692
$ws.Range("G5").Value = 15.281766
$ws.Range("H5").Value = 45.845298
$ws.Range("I5").Value = 0.1817381432449346
$ws.Range("J5").Value = 0.1817381432449346
$ws.Range("M5").Value = 33.83466466666667
$ws.Range("N5").Value = 101.503994
$ws.Range("O5").Value = 0.3204069170612486
$ws.Range("P5").Value = 0.3204069170612486
$ws.Range("Q5").Value = 517.053428124468
$ws.Range("R5").Value = 4653.480853120212
$ws.Range("S5").Value = 0.05823015818954508
$ws.Range("T5").Value = 0.05823015818954508
$ws.Range("I6").Value = 0.2947137116012682
$ws.Range("J6").Value = 0.2947137116012682
$ws.Range("M6").Value = 1.910418
$ws.Range("N6").Value = 5.731254
$ws.Range("O6").Value = 0.01809124304049503
$ws.Range("P6").Value = 0.01809124304049503
$ws.Range("Q6").Value = 47.343046592016
$ws.Range("R6").Value = 426.087419328144
$ws.Range("S6").Value = 0.005331737383944902
$ws.Range("T6").Value = 0.005331737383944902
$ws.Range("I7").Value = 0.2947137116012682
$ws.Range("J7").Value = 0.2947137116012682
$ws.Range("O7").Value = 0.302988173785169
$ws.Range("P7").Value = 0.302988173785169
$ws.Range("S7").Value = 0.08929476926751724
$ws.Range("T7").Value = 0.08929476926751724
$ws.Range("I8").Value = 0.2947137116012682
$ws.Range("J8").Value = 0.2947137116012682
$ws.Range("M8").Value = 37.858701
$ws.Range("N8").Value = 113.576103
$ws.Range("O8").Value = 0.3585136661130873
$ws.Range("P8").Value = 0.3585136661130873
$ws.Range("Q8").Value = 938.195853135912
$ws.Range("R8").Value = 8443.762678223207
$ws.Range("S8").Value = 0.1056588931999658
$ws.Range("T8").Value = 0.1056588931999658
$ws.Range("I9").Value = 0.2947137116012682
$ws.Range("J9").Value = 0.2947137116012682
$ws.Range("M9").Value = 33.83466466666667
$ws.Range("N9").Value = 101.503994
$ws.Range("O9").Value = 0.3204069170612486
$ws.Range("P9").Value = 0.3204069170612486
$ws.Range("Q9").Value = 838.4741484529761
$ws.Range("R9").Value = 7546.267336076785
$ws.Range("S9").Value = 0.09442831174984029
$ws.Range("T9").Value = 0.09442831174984029
$ws.Range("G10").Value = 18.371237
$ws.Range("H10").Value = 55.113711
$ws.Range("I10").Value = 0.2184796247693259
$ws.Range("J10").Value = 0.2184796247693259
$ws.Range("M10").Value = 1.910418
$ws.Range("N10").Value = 5.731254
$ws.Range("O10").Value = 0.01809124304049503
$ws.Range("P10").Value = 0.01809124304049503
$ws.Range("Q10").Value = 35.096741847066
$ws.Range("R10").Value = 315.870676623594
$ws.Range("S10").Value = 0.003952567991098033
$ws.Range("T10").Value = 0.003952567991098033
$ws.Range("G11").Value = 18.371237
$ws.Range("H11").Value = 55.113711
$ws.Range("I11").Value = 0.2184796247693259
$ws.Range("J11").Value = 0.2184796247693259
$ws.Range("O11").Value = 0.302988173785169
$ws.Range("P11").Value = 0.302988173785169
$ws.Range("Q11").Value = 587.792541079094
$ws.Range("R11").Value = 5290.132869711846
$ws.Range("S11").Value = 0.06619674251812704
$ws.Range("T11").Value = 0.06619674251812704
$ws.Range("G12").Value = 18.371237
$ws.Range("H12").Value = 55.113711
$ws.Range("I12").Value = 0.2184796247693259
$ws.Range("J12").Value = 0.2184796247693259
$ws.Range("M12").Value = 37.858701
$ws.Range("N12").Value = 113.576103
$ws.Range("O12").Value = 0.3585136661130873
$ws.Range("P12").Value = 0.3585136661130873
$ws.Range("Q12").Value = 695.511168583137
$ws.Range("R12").Value = 6259.600517248233
$ws.Range("S12").Value = 0.07832793124706272
$ws.Range("T12").Value = 0.07832793124706272
$ws.Range("G13").Value = 18.371237
$ws.Range("H13").Value = 55.113711
$ws.Range("I13").Value = 0.2184796247693259
$ws.Range("J13").Value = 0.2184796247693259
$ws.Range("M13").Value = 33.83466466666667
$ws.Range("N13").Value = 101.503994
$ws.Range("O13").Value = 0.3204069170612486
$ws.Range("P13").Value = 0.3204069170612486
$ws.Range("Q13").Value = 621.5846434068594
$ws.Range("R13").Value = 5594.261790661734
$ws.Range("S13").Value = 0.07000238301303813
$ws.Range("T13").Value = 0.07000238301303813
$ws.Range("G14").Value = 25.652214
$ws.Range("H14").Value = 76.956642
$ws.Range("I14").Value = 0.3050685203844711
$ws.Range("J14").Value = 0.3050685203844711
$ws.Range("M14").Value = 1.910418
$ws.Range("N14").Value = 5.731254
$ws.Range("O14").Value = 0.01809124304049503
$ws.Range("P14").Value = 0.01809124304049503
$ws.Range("Q14").Value = 49.006451365452
$ws.Range("R14").Value = 441.058062289068
$ws.Range("S14").Value = 0.005519068746279679
$ws.Range("T14").Value = 0.005519068746279679
$ws.Range("G15").Value = 25.652214
$ws.Range("H15").Value = 76.956642
$ws.Range("I15").Value = 0.3050685203844711
$ws.Range("J15").Value = 0.3050685203844711
$ws.Range("O15").Value = 0.302988173785169
$ws.Range("P15").Value = 0.302988173785169
$ws.Range("Q15").Value = 820.749307810068
$ws.Range("R15").Value = 7386.743770290613
$ws.Range("S15").Value = 0.09243215387063451
$ws.Range("T15").Value = 0.09243215387063451
$ws.Range("G16").Value = 25.652214
$ws.Range("H16").Value = 76.956642
$ws.Range("I16").Value = 0.3050685203844711
$ws.Range("J16").Value = 0.3050685203844711
$ws.Range("M16").Value = 37.858701
$ws.Range("N16").Value = 113.576103
$ws.Range("O16").Value = 0.3585136661130873
$ws.Range("P16").Value = 0.3585136661130873
$ws.Range("Q16").Value = 971.1594998140139
$ws.Range("R16").Value = 8740.435498326126
$ws.Range("S16").Value = 0.1093712336587319
$ws.Range("T16").Value = 0.1093712336587319
$ws.Range("G17").Value = 25.652214
$ws.Range("H17").Value = 76.956642
$ws.Range("I17").Value = 0.3050685203844711
$ws.Range("J17").Value = 0.3050685203844711
$ws.Range("M17").Value = 33.83466466666667
$ws.Range("N17").Value = 101.503994
$ws.Range("O17").Value = 0.3204069170612486
$ws.Range("P17").Value = 0.3204069170612486
$ws.Range("Q17").Value = 867.9340586475721
$ws.Range("R17").Value = 7811.406527828149
$ws.Range("S17").Value = 0.09774606410882507
$ws.Range("T17").Value = 0.09774606410882507
